# Edit applied to Saldo.xlsx "Export" sheet.
# Working from the bottom of the affected block upward so row numbers
# of not-yet-processed rows stay stable while we delete/insert.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the trailing trio of rows that are dropped entirely:
#   row 9 -> 005654122 / ELIANE   / 3748.74
#   row 8 -> 004313254 / GUSTAVO  / 4292   (this account is kept, but moved
#                                           up to replace the SARA row below,
#                                           so this duplicate copy goes away)
#   row 7 -> 005018038 / ELAINE   / 5562.91
$ws.Rows.Item(9).Delete()
$ws.Rows.Item(8).Delete()
$ws.Rows.Item(7).Delete()

# Row 6 (004511696 / KRYSCIA) keeps its account/name, only the balance changes.
$ws.Cells.Item(6, 3).Value = 1450

# Row 5 used to be 008032257 / SARA / 18922.37; it becomes the GUSTAVO row
# that used to sit further down (004313254 / GUSTAVO / 4292).
# Account numbers are zero-padded strings, not numbers, so force text
# formatting before assigning them (otherwise Excel would coerce them to
# numbers and the leading zeros would be lost).
$ws.Cells.Item(5, 1).NumberFormat = "@"
$ws.Cells.Item(5, 1).Value = "004313254"
$ws.Cells.Item(5, 2).Value = "GUSTAVO"
$ws.Cells.Item(5, 3).Value = 4292

# Row 3 (004556150 / MARINA / 50104.17) is removed outright.
$ws.Rows.Item(3).Delete()

# Insert two brand-new rows right above the existing 004368468 / AHMAD row
# (currently row 2).
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(2).Insert()

$ws.Cells.Item(2, 1).NumberFormat = "@"
$ws.Cells.Item(2, 1).Value = "001882235"
$ws.Cells.Item(2, 2).Value = "LAGO"
$ws.Cells.Item(2, 3).Value = 278323.91

$ws.Cells.Item(3, 1).NumberFormat = "@"
$ws.Cells.Item(3, 1).Value = "008028807"
$ws.Cells.Item(3, 2).Value = "RAFAEL"
$ws.Cells.Item(3, 3).Value = 100000
